# Add team record columns (Wins / Losses / Ties) to the right of the
# existing data, matching the header styling already used for row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, thin border, centered alignment)
# from the last existing header cell (AB1) onto the three new header cells.
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header labels
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Team record values for every player row (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 29).Value = 70   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 92   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
